$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 252, shifting existing rows 252:267 down to 253:268
# (matches the weekly price-update pattern: a new week's reading is
# prepended and the rest of the "Ajo" series shifts down by one row).
$ws.Rows.Item(252).Insert()

# Populate the new row 252 with this week's reading. The "dimension"
# columns (Mercado/Region/Categoria/etc.) repeat the same constant values
# used throughout this Ajo block.
$ws.Range("A252").Value = 4
$ws.Range("B252").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C252").Value = "Los Lagos"
$ws.Range("D252").Value = 44706
$ws.Range("E252").Value = 10
$ws.Range("F252").Value = 100112003
$ws.Range("G252").Value = "Ajo"
$ws.Range("H252").Value = "Chino"
$ws.Range("I252").Value = "Primera"
$ws.Range("J252").Value = 15
$ws.Range("K252").Value = 22000
$ws.Range("L252").Value = 22000
$ws.Range("M252").Value = 22000
$ws.Range("N252").Value = "$/caja 10 kilos"
$ws.Range("O252").Value = "China"
$ws.Range("P252").Value = 2200
$ws.Range("Q252").Value = 10
$ws.Range("R252").Value = "Hortaliza"
